$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 0.8156942129135132
$ws.Range("B1").Value = 1.181926608085632
$ws.Range("C1").Value = 2.309945583343506
$ws.Range("D1").Value = 3.946348667144775
$ws.Range("E1").Value = 1.920693159103394
